# floodchannel.xlsx edit script
# - trials 1D-2D coupling, small change to channel WH
# - cleaned up plot graph for P and Qs
# - linked flood check then no erosion

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core input changes (left block, column C / middle block, column G) ---

# Row 7: cwh inputs
$ws.Range("C7").Value = 2
$ws.Range("G7").Value = 2

# Row 8: cd inputs; G8 becomes a plain value instead of "=C8"
$ws.Range("C8").Value = 3
$ws.Range("G8").Value = 0.8

# Row 11: cw inputs
$ws.Range("C11").Value = 0.07
$ws.Range("G11").Value = 0.05

# Row 12: adj inputs
$ws.Range("C12").Value = 0.03
$ws.Range("G12").Value = 0

# Row 16: G16 now clamps the water height drop to a minimum of zero
$ws.Range("G16").Formula = "=MAX(G15,0)"

# --- New lower rows: "if < 0" / "hmx/fc" labels and follow-on calcs ---

$ws.Range("F21").Value = "if < 0"
$ws.Range("G21").Formula = "=G11*G9/G7"
$ws.Range("H21").Value = "hmx/fc"

$ws.Range("C22").Formula = "=C12+C11/(C7/C9)"
$ws.Range("G22").Formula = "=G21*2"

$ws.Range("C23").Formula = "=C22*2"
$ws.Range("G24").Formula = "=G12+G11*(G9/G7)"

# remove the old standalone SQRT check and the now-empty G20/G23/G24 style cells
$ws.Range("G20").Clear()

# de-bold rows 23-26 in columns B (regular) and de-bold (to italic) column C
$ws.Range("B23:B26").Font.Bold = $false
$ws.Range("C23:C26").Font.Bold = $false

# clear the old D29 mass-check formula, keep it blank but bold-styled like its neighbours
$ws.Range("D29").ClearContents()
$ws.Range("D29").Font.Bold = $true

# --- View state: scroll down a bit and move the active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()
